# Insert a new data row at row 86, pushing existing rows 86-130 down to 87-131.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("86:86").Insert()

# Populate the newly inserted row 86 with the new record's data.
$ws.Range("A86").Value = 11
$ws.Range("B86").Value = "Vega Monumental Concepción"
$ws.Range("C86").Value = "Bíobío"
$ws.Range("D86").Value = 45146
$ws.Range("E86").Value = 8
$ws.Range("F86").Value = 100112037
$ws.Range("G86").Value = "Cebollín"
$ws.Range("H86").Value = "Sin especificar"
$ws.Range("I86").Value = "Primera"
$ws.Range("J86").Value = 100
$ws.Range("K86").Value = 5000
$ws.Range("L86").Value = 5500
$ws.Range("M86").Value = 5250
$ws.Range("N86").Value = "$/paquete 36 unidades"
$ws.Range("O86").Value = "Región Metropolitana"
$ws.Range("P86").Value = 146
$ws.Range("Q86").Value = 36
$ws.Range("R86").Value = "Hortaliza"
